# Add new book entry "You Never Forget Your First" as row 74 on the
# "Completed" sheet (the active sheet of the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title / Author
$ws.Range("A74").Value = "You Never Forget Your First"
$ws.Range("B74").Value = "Alexis Coe"

# Start Date - copy the date formatting from the cell above (row 73)
# first, then set the value so the existing date style (s="1") is
# reused instead of a brand-new number format being created.
[void]$ws.Range("C73").Copy()
$ws.Range("C74").PasteSpecial(-4122)
$ws.Range("C74").Value = (Get-Date -Year 2020 -Month 5 -Day 7 -Hour 0 -Minute 0 -Second 0)

# Finish Date
[void]$ws.Range("D73").Copy()
$ws.Range("D74").PasteSpecial(-4122)
$ws.Range("D74").Value = (Get-Date -Year 2020 -Month 5 -Day 8 -Hour 0 -Minute 0 -Second 0)

# Tags
$ws.Range("E74").Value = "george washington;history;biography;president;revolutionary war;american history"

# Type (reuses the existing "Audio" shared string)
$ws.Range("F74").Value = "Audio"

# Length
$ws.Range("G74").Value = "6 Hours 42 Mins"

# Move the selection down to the next empty row, like Excel does after
# typing a new row of data.
[void]$ws.Range("A75").Select()
